$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 393, shifting existing rows 393:423 down to 394:424
$ws.Rows("393:393").Insert()

# Populate the newly inserted row 393 with the new weekly record
$ws.Range("A393").Value = 10
$ws.Range("B393").Value = "Vega Modelo de Temuco"
$ws.Range("C393").Value = "La Araucanía"
$ws.Range("D393").Value = (Get-Date -Year 2022 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E393").Value = 9
$ws.Range("F393").Value = "Fruta"
$ws.Range("G393").Value = 100108
$ws.Range("H393").Value = "Tropicales y subtropicales"
$ws.Range("I393").Value = 100108002
$ws.Range("J393").Value = "Mango"
$ws.Range("K393").Value = "Sin especificar"
$ws.Range("L393").Value = "Primera"
$ws.Range("M393").Value = 350
$ws.Range("N393").Value = 10000
$ws.Range("O393").Value = 10000
$ws.Range("P393").Value = 10000
$ws.Range("Q393").Value = "`$/bandeja 4 kilos"
$ws.Range("R393").Value = "Brasil"
$ws.Range("S393").Value = 2500
$ws.Range("T393").Value = 4
